# Convert the bsecode (column D) for existing rows 727-746 from text to numeric,
# matching how the rest of the "day" sheet stores bsecode values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

$ws.Cells.Item(727, 4).Value = 508869
$ws.Cells.Item(728, 4).Value = 500034
$ws.Cells.Item(729, 4).Value = 540005
$ws.Cells.Item(730, 4).Value = 500830
$ws.Cells.Item(731, 4).Value = 500331
$ws.Cells.Item(732, 4).Value = 503806
$ws.Cells.Item(733, 4).Value = 500410
$ws.Cells.Item(734, 4).Value = 524715
$ws.Cells.Item(735, 4).Value = 532755
$ws.Cells.Item(736, 4).Value = 500043
$ws.Cells.Item(737, 4).Value = 500425
$ws.Cells.Item(738, 4).Value = 500295
$ws.Cells.Item(739, 4).Value = 532555
$ws.Cells.Item(740, 4).Value = 500104
$ws.Cells.Item(741, 4).Value = 500547
$ws.Cells.Item(742, 4).Value = 532155
$ws.Cells.Item(743, 4).Value = 530965
$ws.Cells.Item(744, 4).Value = 500470
$ws.Cells.Item(745, 4).Value = 532483
$ws.Cells.Item(746, 4).Value = 539437

# Append the newly scraped rows (747-759) for the 16/10/2024 11:36:28 batch.
# bsecode (column D) keeps its original text/string representation, matching
# the source data feed -- use a leading apostrophe so Excel stores it as text
# rather than auto-converting the numeric-looking code to a number.

$ws.Cells.Item(747, 1).Value = 1
$ws.Cells.Item(747, 2).Value = "NAVINFLUOR"
$ws.Cells.Item(747, 3).Value = "Navin Fluorine International Limited"
$ws.Cells.Item(747, 4).Value = "'532504"
$ws.Cells.Item(747, 5).Value = 3.03
$ws.Cells.Item(747, 6).Value = 3461.25
$ws.Cells.Item(747, 7).Value = 286792
$ws.Cells.Item(747, 8).Value = "day"
$ws.Cells.Item(747, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(748, 1).Value = 2
$ws.Cells.Item(748, 2).Value = "TATACOMM"
$ws.Cells.Item(748, 3).Value = "Tata Communications Limited"
$ws.Cells.Item(748, 4).Value = "'500483"
$ws.Cells.Item(748, 5).Value = -0.46
$ws.Cells.Item(748, 6).Value = 1922.45
$ws.Cells.Item(748, 7).Value = 210463
$ws.Cells.Item(748, 8).Value = "day"
$ws.Cells.Item(748, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(749, 1).Value = 3
$ws.Cells.Item(749, 2).Value = "AUROPHARMA"
$ws.Cells.Item(749, 3).Value = "Aurobindo Pharma Limited"
$ws.Cells.Item(749, 4).Value = "'524804"
$ws.Cells.Item(749, 5).Value = 0.14
$ws.Cells.Item(749, 6).Value = 1482.05
$ws.Cells.Item(749, 7).Value = 606485
$ws.Cells.Item(749, 8).Value = "day"
$ws.Cells.Item(749, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(750, 1).Value = 4
$ws.Cells.Item(750, 2).Value = "CANFINHOME"
$ws.Cells.Item(750, 3).Value = "Can Fin Homes Limited"
$ws.Cells.Item(750, 4).Value = "'511196"
$ws.Cells.Item(750, 5).Value = -0.95
$ws.Cells.Item(750, 6).Value = 882.35
$ws.Cells.Item(750, 7).Value = 194324
$ws.Cells.Item(750, 8).Value = "day"
$ws.Cells.Item(750, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(751, 1).Value = 5
$ws.Cells.Item(751, 2).Value = "ICICIPRULI"
$ws.Cells.Item(751, 3).Value = "Icici Prudential Life Insurance Company Limited"
$ws.Cells.Item(751, 4).Value = "'540133"
$ws.Cells.Item(751, 5).Value = 2
$ws.Cells.Item(751, 6).Value = 748.45
$ws.Cells.Item(751, 7).Value = 1139316
$ws.Cells.Item(751, 8).Value = "day"
$ws.Cells.Item(751, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(752, 1).Value = 6
$ws.Cells.Item(752, 2).Value = "UPL"
$ws.Cells.Item(752, 3).Value = "Upl Limited"
$ws.Cells.Item(752, 4).Value = "'512070"
$ws.Cells.Item(752, 5).Value = -0.91
$ws.Cells.Item(752, 6).Value = 568.85
$ws.Cells.Item(752, 7).Value = 1303742
$ws.Cells.Item(752, 8).Value = "day"
$ws.Cells.Item(752, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(753, 1).Value = 7
$ws.Cells.Item(753, 2).Value = "RECLTD"
$ws.Cells.Item(753, 3).Value = "Rural Electrification Corporation Limited"
$ws.Cells.Item(753, 4).Value = "'532955"
$ws.Cells.Item(753, 5).Value = 0.88
$ws.Cells.Item(753, 6).Value = 552.1
$ws.Cells.Item(753, 7).Value = 5230745
$ws.Cells.Item(753, 8).Value = "day"
$ws.Cells.Item(753, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(754, 1).Value = 8
$ws.Cells.Item(754, 2).Value = "IGL"
$ws.Cells.Item(754, 3).Value = "Indraprastha Gas Limited"
$ws.Cells.Item(754, 4).Value = "'532514"
$ws.Cells.Item(754, 5).Value = -1.16
$ws.Cells.Item(754, 6).Value = 518.55
$ws.Cells.Item(754, 7).Value = 1146070
$ws.Cells.Item(754, 8).Value = "day"
$ws.Cells.Item(754, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(755, 1).Value = 9
$ws.Cells.Item(755, 2).Value = "PFC"
$ws.Cells.Item(755, 3).Value = "Power Finance Corporation Limited"
$ws.Cells.Item(755, 4).Value = "'532810"
$ws.Cells.Item(755, 5).Value = 0.51
$ws.Cells.Item(755, 6).Value = 479.2
$ws.Cells.Item(755, 7).Value = 5930214
$ws.Cells.Item(755, 8).Value = "day"
$ws.Cells.Item(755, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(756, 1).Value = 10
$ws.Cells.Item(756, 2).Value = "HINDCOPPER"
$ws.Cells.Item(756, 3).Value = "Hindustan Copper Limited"
$ws.Cells.Item(756, 4).Value = "'513599"
$ws.Cells.Item(756, 5).Value = 0.31
$ws.Cells.Item(756, 6).Value = 321.9
$ws.Cells.Item(756, 7).Value = 3987022
$ws.Cells.Item(756, 8).Value = "day"
$ws.Cells.Item(756, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(757, 1).Value = 11
$ws.Cells.Item(757, 2).Value = "NMDC"
$ws.Cells.Item(757, 3).Value = "Nmdc Limited"
$ws.Cells.Item(757, 4).Value = "'526371"
$ws.Cells.Item(757, 5).Value = -0.29
$ws.Cells.Item(757, 6).Value = 231.4
$ws.Cells.Item(757, 7).Value = 9649158
$ws.Cells.Item(757, 8).Value = "day"
$ws.Cells.Item(757, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(758, 1).Value = 12
$ws.Cells.Item(758, 2).Value = "FEDERALBNK"
$ws.Cells.Item(758, 3).Value = "The Federal Bank  Limited"
$ws.Cells.Item(758, 4).Value = "'500469"
$ws.Cells.Item(758, 5).Value = -2.12
$ws.Cells.Item(758, 6).Value = 194.4
$ws.Cells.Item(758, 7).Value = 8161660
$ws.Cells.Item(758, 8).Value = "day"
$ws.Cells.Item(758, 9).Value = "16/10/2024 11:36:28"

$ws.Cells.Item(759, 1).Value = 13
$ws.Cells.Item(759, 2).Value = "LTF"
$ws.Cells.Item(759, 3).Value = "L&T Finance Ltd"
$ws.Cells.Item(759, 4).Value = "'533519"
$ws.Cells.Item(759, 5).Value = -0.73
$ws.Cells.Item(759, 6).Value = 166.83
$ws.Cells.Item(759, 7).Value = 7759481
$ws.Cells.Item(759, 8).Value = "day"
$ws.Cells.Item(759, 9).Value = "16/10/2024 11:36:28"

